$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: new date entry (matching style of A1/A6 - bold, d-mmm-yy format)
# NOTE: format/font must be applied BEFORE the value is written, otherwise
# Excel auto-assigns a transient "m/d/yyyy" number format to the date value
# that lingers as an unused numFmt entry in styles.xml.
$ws.Range("A8").NumberFormat = "d-mmm-yy"
$ws.Range("A8").Font.Bold = $true
$d = Get-Date -Year 2022 -Month 5 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Range("A8").Value = $d

# Entries are written in the order the shared strings were originally authored
# (row 10's comment #2 first, then row 11's comment #3, then row 9's comment #1,
# then row 12's comment #4), to match the shared-string table ordering.

# Row 10: comment #2 (carries a numFmt 15 / non-bold style, like a leftover date format applied to text)
$ws.Range("A10").Value = "2. The old data set only has sales through the first part of 2016, need to drop year==2016, so that we rely on the new data set for 2016 sales"
$ws.Range("A10").NumberFormat = "d-mmm-yy"
$ws.Range("A10").Font.Bold = $false

# Row 11: comment #3
$ws.Range("A11").Value = "3. the old data set (9516) contains some variables that the new one (1621) does not. Drop variable from 9516 if it does not appear in 1621"

# Row 9: comment #1
$ws.Range("A9").Value = "1. Started data merge process of ohiohousesales_9516_cleaned.dta (old) and ohiohousesales_1621_cleaned.dta (new). Goal is to append the two datasets together "

# Row 12: comment #4
$ws.Range("A12").Value = "4. Variable names differed between datasets. I created a mapping table "

# Update selection to match final state (I7 was selected when saved)
$ws.Range("I7").Select() | Out-Null

$wb.Save()
